# Post adding of adjustments functionality
# Updates the "cash" / "accounts_receivable" (renamed "ar") rows with new
# figures, and appends new balance-sheet line items (subtotal/total rows for
# assets, liabilities and equity) below the existing two data rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# All "numbers" in this sheet are stored as text (inline strings), not
# numeric values - so we force Text formatting before writing, then strip
# the formatting override back off so the cells end up with no explicit
# style, matching a plain inline-string cell.
$dataRange = $ws.Range("A2:G13")
$dataRange.NumberFormat = "@"

# ---- Row 2: Cash & Equivalents (unchanged label/id, new figures) ----
$ws.Range("D2").Value = "50"
$ws.Range("E2").Value = "55"
$ws.Range("F2").Value = "61"
$ws.Range("G2").Value = "67"

# ---- Row 3: Accounts Receivable -> id renamed to "ar", new figures ----
$ws.Range("B3").Value = "ar"
$ws.Range("D3").Value = "100"
$ws.Range("E3").Value = "112"
$ws.Range("F3").Value = "125"
# G3 stays "140" - no change

# ---- Row 4: Total Current Assets (new) ----
$ws.Range("A4").Value = "  Total Current Assets"
$ws.Range("B4").Value = "current_assets_subtotal"
$ws.Range("C4").Value = "150"
$ws.Range("D4").Value = "150"
$ws.Range("E4").Value = "167"
$ws.Range("F4").Value = "186"
$ws.Range("G4").Value = "207"

# ---- Row 5: Property, Plant & Equipment (new) ----
$ws.Range("A5").Value = "  Property, Plant & Equipment"
$ws.Range("B5").Value = "ppe"
$ws.Range("C5").Value = "300"
$ws.Range("D5").Value = "300"
$ws.Range("E5").Value = "315"
$ws.Range("F5").Value = "331"
$ws.Range("G5").Value = "347"

# ---- Row 6: Total Assets (new) ----
$ws.Range("A6").Value = "  Total Assets"
$ws.Range("B6").Value = "total_assets"
$ws.Range("C6").Value = "450"
$ws.Range("D6").Value = "450"
$ws.Range("E6").Value = "482"
$ws.Range("F6").Value = "517"
$ws.Range("G6").Value = "554"

# ---- Row 7: Accounts Payable (new) ----
$ws.Range("A7").Value = "  Accounts Payable"
$ws.Range("B7").Value = "ap"
$ws.Range("C7").Value = "80"
$ws.Range("D7").Value = "80"
$ws.Range("E7").Value = "86"
$ws.Range("F7").Value = "93"
$ws.Range("G7").Value = "101"

# ---- Row 8: Long-Term Debt (new) ----
$ws.Range("A8").Value = "  Long-Term Debt"
$ws.Range("B8").Value = "debt"
$ws.Range("C8").Value = "150"
$ws.Range("D8").Value = "150"
$ws.Range("E8").Value = "153"
$ws.Range("F8").Value = "156"
$ws.Range("G8").Value = "159"

# ---- Row 9: Total Liabilities (new) ----
$ws.Range("A9").Value = "  Total Liabilities"
$ws.Range("B9").Value = "total_liabilities"
$ws.Range("C9").Value = "230"
$ws.Range("D9").Value = "230"
$ws.Range("E9").Value = "239"
$ws.Range("F9").Value = "249"
$ws.Range("G9").Value = "260"

# ---- Row 10: Common Stock (new) ----
$ws.Range("A10").Value = "  Common Stock"
$ws.Range("B10").Value = "common_stock"
$ws.Range("C10").Value = "100"
$ws.Range("D10").Value = "100"
$ws.Range("E10").Value = "100"
$ws.Range("F10").Value = "100"
$ws.Range("G10").Value = "100"

# ---- Row 11: Retained Earnings (new) ----
$ws.Range("A11").Value = "  Retained Earnings"
$ws.Range("B11").Value = "retained_earnings"
$ws.Range("C11").Value = "1,035"
$ws.Range("D11").Value = "2,180"
$ws.Range("E11").Value = "3,527"
$ws.Range("F11").Value = "5,114"
$ws.Range("G11").Value = "7,039"

# ---- Row 12: Total Equity (new) ----
$ws.Range("A12").Value = "  Total Equity"
$ws.Range("B12").Value = "total_equity"
$ws.Range("C12").Value = "1,135"
$ws.Range("D12").Value = "2,280"
$ws.Range("E12").Value = "3,627"
$ws.Range("F12").Value = "5,214"
$ws.Range("G12").Value = "7,139"

# ---- Row 13: Total Liabilities & Equity (new) ----
$ws.Range("A13").Value = "  Total Liabilities & Equity"
$ws.Range("B13").Value = "total_liabs_equity"
$ws.Range("C13").Value = "1,365"
$ws.Range("D13").Value = "2,510"
$ws.Range("E13").Value = "3,867"
$ws.Range("F13").Value = "5,464"
$ws.Range("G13").Value = "7,399"

# Strip the temporary Text number-format override back off so the cells
# don't carry an explicit style (matches the target's plain inline strings).
$dataRange.ClearFormats()
